# pedalboard-display-bom.xlsx update
# - BoM sheet: J1 connector value "LED-IN" -> "JST PH 3"
#              J3 connector value "OLED"   -> "JST PH 4"
#              OLED (U1 U2) row gets a real Value ("OLED I2C") plus a
#              Datasheet + Supplier link (previously empty / placeholder cells)
# - DNF sheet: J2 connector value "LED-OUT" -> "JST PH 3" (matches BoM rename)
#              Resistor (R2 R3) row gets a Datasheet + Supplier link
#              Datasheet/Supplier columns widened to fit the new long URLs
#              and the resistor row is made taller to match

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# BoM sheet
# ---------------------------------------------------------------------
$bom = $wb.Worksheets.Item("BoM")

# J1 connector - Value column
$bom.Range("E11").Value = "JST PH 3"

# J3 connector - Value column
$bom.Range("E12").Value = "JST PH 4"

# OLED (U1 U2) row - Value / Datasheet / Supplier
$bom.Range("E14").Value = "OLED I2C"
$bom.Range("I14").Value = "img/GME12812.pdf"
$bom.Range("J14").Value = "https://www.aliexpress.com/item/1005005253671396.html"

# These three cells used to carry the "empty/placeholder" row-banding fill;
# now that they hold real data they pick up the "has data" banding fill,
# same as the other populated rows (e.g. row 10 / row 12).
$bom.Range("E14").Interior.Color = 12447999   # RGB(189,240,255) -> FFFFF0BD fill
$bom.Range("I14").Interior.Color = 12447999   # FFFFF0BD fill
$bom.Range("J14").Interior.Color = 16777200   # RGB(255,255,240) -> FFF0FFFF fill

# ---------------------------------------------------------------------
# DNF sheet
# ---------------------------------------------------------------------
$dnf = $wb.Worksheets.Item("DNF")

# J2 connector - Value column (renamed to match the BoM sheet)
$dnf.Range("E9").Value = "JST PH 3"

# Resistor (R2 R3) row - Datasheet / Supplier
$dnf.Range("I10").Value = "https://www.seielect.com/catalog/sei-rmcf_rmcp.pdf"
$dnf.Range("J10").Value = "https://www.digikey.ch/en/products/detail/stackpole-electronics-inc/RMCF0805FT3K30/1760325"

$dnf.Range("I10").Interior.Color = 12447999   # FFFFF0BD fill
$dnf.Range("J10").Interior.Color = 16777200   # FFF0FFFF fill

# Row is now taller to comfortably show the wrapped URLs
$dnf.Rows.Item(10).RowHeight = 30

# Datasheet / Supplier columns widened for the long URLs
$dnf.Columns.Item(9).ColumnWidth = 54.8
$dnf.Columns.Item(10).ColumnWidth = 59.8
